# Applies the "Updated cryptos list" price/volume refresh described in the
# commit diff. Text-ish numeric-looking values (prices like "1.001",
# "0.8960", "27.365.18") must stay plain text, matching the workbook's
# original inlineStr cells, so each write forces Text number format first
# and restores the default "Normal" style afterwards (no lingering @ format).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.365.18"
Set-TextValue "E2" "  +1.40%  "
Set-TextValue "D3" "1.827.94"
Set-TextValue "E3" "  +0.30%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "313.93"
Set-TextValue "E5" "  +1.07%  "
Set-TextValue "E6" "  +0.01%  "
Set-TextValue "D7" "0.4474"
Set-TextValue "E7" "  +5.25%  "
Set-TextValue "D8" "0.3769"
Set-TextValue "E8" "  +3.11%  "
Set-TextValue "D9" "0.07524"
Set-TextValue "E9" "  +4.13%  "
Set-TextValue "D10" "0.8960"
Set-TextValue "E10" "  +6.50%  "
Set-TextValue "D11" "21.07"
Set-TextValue "E11" "  +2.58%  "
Set-TextValue "D12" "1.808.04"
Set-TextValue "E12" "  -1.04%  "
Set-TextValue "D13" "6.741"
Set-TextValue "E13" "  +1.42%  "
Set-TextValue "D14" "94.56"
Set-TextValue "E14" "  +5.77%  "
Set-TextValue "D15" "5.406"
Set-TextValue "E15" "  +2.50%  "
Set-TextValue "D16" "0.07116"
Set-TextValue "E16" "  +0.74%  "
Set-TextValue "D17" "1.002"
Set-TextValue "E17" "  -0.05%  "
Set-TextValue "D18" "0.000008832"
Set-TextValue "D19" "1.001"
Set-TextValue "E19" "  -0.03%  "
Set-TextValue "D20" "15.24"
Set-TextValue "E20" "  +2.67%  "
Set-TextValue "D21" "27.389.28"
Set-TextValue "E21" "  +0.79%  "
Set-TextValue "E22" "  +3.17%  "
Set-TextValue "E23" "  +1.53%  "
Set-TextValue "D24" "2.003"
Set-TextValue "E24" "  +1.41%  "
Set-TextValue "D25" "2.501"
Set-TextValue "E25" "  +12.76%  "
Set-TextValue "D26" "151.59"
Set-TextValue "E26" "  +0.46%  "
Set-TextValue "D27" "18.59"
Set-TextValue "E27" "  +2.38%  "
Set-TextValue "D28" "5.366"
Set-TextValue "E28" "  +3.07%  "
Set-TextValue "D29" "118.09"
Set-TextValue "E29" "  +1.24%  "
Set-TextValue "D30" "0.08848"
Set-TextValue "E30" "  +1.54%  "
Set-TextValue "D31" "0.7816"
Set-TextValue "E31" "  +6.52%  "
Set-TextValue "E32" "  +2.49%  "
Set-TextValue "D33" "4.528"
Set-TextValue "E33" "  +2.64%  "
Set-TextValue "D34" "2.892"
Set-TextValue "E34" "  -0.34%  "
Set-TextValue "D35" "1.000"
Set-TextValue "E35" "  +0.00%  "
Set-TextValue "D36" "1.109"
Set-TextValue "E36" "  +1.92%  "
Set-TextValue "D37" "0.01991"
Set-TextValue "E37" "  +2.81%  "
Set-TextValue "D38" "0.05336"
Set-TextValue "E38" "  +2.47%  "
Set-TextValue "D39" "7.387"
Set-TextValue "E39" "  +2.52%  "
Set-TextValue "D40" "0.5330"
Set-TextValue "E40" "  +4.42%  "
Set-TextValue "D43" "2.279"
Set-TextValue "E43" "  +16.82%  "
Set-TextValue "D44" "8.820"
Set-TextValue "E44" "  +3.59%  "
Set-TextValue "D45" "0.5180"
Set-TextValue "E45" "  +9.67%  "
Set-TextValue "D46" "10.73"
Set-TextValue "E46" "  +2.28%  "
Set-TextValue "D47" "106.29"
Set-TextValue "E47" "  +0.62%  "
Set-TextValue "D48" "1.709"
Set-TextValue "E48" "  +3.81%  "
Set-TextValue "D49" "1.000"
Set-TextValue "E49" "  +0.06%  "
Set-TextValue "D50" "0.06375"
Set-TextValue "E50" "  +0.89%  "
Set-TextValue "D51" "64.54"
Set-TextValue "E51" "  +3.45%  "

# Rows 41 and 42: Algorand and MXToken swap positions, each with refreshed
# price/volume figures.
Set-TextValue "B41" "MXToken"
Set-TextValue "C41" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.878"
Set-TextValue "E41" "  +0.40%  "

Set-TextValue "B42" "Algorand"
Set-TextValue "C42" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D42" "0.1731"
Set-TextValue "E42" "  +2.99%  "
